$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text updates (rich-text cells: patch only the changed substring) ----
$ws.Range("A8").Characters(21,2).Text = "27"
$ws.Range("C9").Characters(27,9).Text = "7/1/2024"
$ws.Range("C9").Characters(47,9).Text = "7/7/2024"

# ---- Cells whose type/style changes (number <-> shared-text) ----
# Strategy: set the literal value first (apostrophe-prefix forces text for
# numeric-looking strings), then copy ONLY the number-format/style from an
# unaffected donor cell that already carries the desired style so the cell
# ends up with the same style index as the target.

$ws.Range("C15").Value = "'0"
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "'***.*"
$ws.Range("D17").Value = "'0"
$ws.Range("E17").Value = "'***.*"
$ws.Range("C18").Value = 2
$ws.Range("D22").Value = "'0"
$ws.Range("E22").Value = "'***.*"
$ws.Range("C27").Value = "'0"
$ws.Range("C28").Value = "'0"
$ws.Range("D28").Value = "'0"
$ws.Range("E28").Value = "'***.*"

# Donor cells for format copy (unaffected by this edit, already carry target style)
$donorText0     = $ws.Range("C14")   # style 14, shared text "0"
$donorTextStar  = $ws.Range("E14")   # style 14, shared text "***.*"
$donorNumber16  = $ws.Range("C19")   # style 16, plain number

$donorText0.Copy()
$ws.Range("C15").PasteSpecial(-4122)
$donorText0.Copy()
$ws.Range("D15").PasteSpecial(-4122)
$donorTextStar.Copy()
$ws.Range("E15").PasteSpecial(-4122)
$donorText0.Copy()
$ws.Range("D17").PasteSpecial(-4122)
$donorTextStar.Copy()
$ws.Range("E17").PasteSpecial(-4122)
$donorNumber16.Copy()
$ws.Range("C18").PasteSpecial(-4122)
$donorText0.Copy()
$ws.Range("D22").PasteSpecial(-4122)
$donorTextStar.Copy()
$ws.Range("E22").PasteSpecial(-4122)
$donorText0.Copy()
$ws.Range("C27").PasteSpecial(-4122)
$donorText0.Copy()
$ws.Range("C28").PasteSpecial(-4122)
$donorText0.Copy()
$ws.Range("D28").PasteSpecial(-4122)
$donorTextStar.Copy()
$ws.Range("E28").PasteSpecial(-4122)

# ---- Plain value updates (style/type unchanged) ----
$ws.Range("N15").Value = 100
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 133.333333333333
$ws.Range("I16").Value = 42
$ws.Range("J16").Value = 36
$ws.Range("K16").Value = 16.666666666666
$ws.Range("L16").Value = 13.513513513513
$ws.Range("M16").Value = -14.285714285714
$ws.Range("N16").Value = -88.135593220339
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 66.666666666666
$ws.Range("I17").Value = 51
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 70
$ws.Range("N17").Value = -22.727272727272
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -30
$ws.Range("I18").Value = 47
$ws.Range("J18").Value = 65
$ws.Range("K18").Value = -27.692307692307
$ws.Range("L18").Value = -14.545454545454
$ws.Range("M18").Value = -28.787878787878
$ws.Range("N18").Value = -93.128654970760
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = 40.625
$ws.Range("I19").Value = 229
$ws.Range("J19").Value = 248
$ws.Range("K19").Value = -7.661290322580
$ws.Range("L19").Value = -12.260536398467
$ws.Range("M19").Value = 30.113636363636
$ws.Range("N19").Value = -54.016064257028
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = -88.888888888888
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = -31.578947368421
$ws.Range("I20").Value = 88
$ws.Range("J20").Value = 72
$ws.Range("K20").Value = 22.222222222222
$ws.Range("L20").Value = 91.304347826087
$ws.Range("M20").Value = 66.037735849056
$ws.Range("N20").Value = -95.070028011204
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -45
$ws.Range("F21").Value = 83
$ws.Range("G21").Value = 71
$ws.Range("H21").Value = 16.901408450704
$ws.Range("I21").Value = 463
$ws.Range("J21").Value = 475
$ws.Range("K21").Value = -2.526315789473
$ws.Range("L21").Value = 0.871459694989
$ws.Range("M21").Value = 23.138297872340
$ws.Range("N21").Value = -86.350235849056
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = -25
$ws.Range("L22").Value = 23.529411764705
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = 30.434782608695
$ws.Range("F24").Value = 149
$ws.Range("G24").Value = 110
$ws.Range("H24").Value = 35.454545454545
$ws.Range("I24").Value = 850
$ws.Range("J24").Value = 808
$ws.Range("K24").Value = 5.198019801980
$ws.Range("L24").Value = -7.407407407407
$ws.Range("M24").Value = 74.537987679671
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = 18.75
$ws.Range("F25").Value = 97
$ws.Range("G25").Value = 79
$ws.Range("H25").Value = 22.784810126582
$ws.Range("I25").Value = 611
$ws.Range("J25").Value = 580
$ws.Range("K25").Value = 5.344827586206
$ws.Range("L25").Value = -6.717557251908
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 66.666666666666
$ws.Range("F26").Value = 22
$ws.Range("G26").Value = 13
$ws.Range("H26").Value = 69.230769230769
$ws.Range("I26").Value = 147
$ws.Range("J26").Value = 122
$ws.Range("K26").Value = 20.491803278688
$ws.Range("L26").Value = 53.125
$ws.Range("M26").Value = 28.947368421052
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 10
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = 42.857142857142
$ws.Range("L27").Value = -9.090909090909
$ws.Range("I28").Value = 19
$ws.Range("K28").Value = 26.666666666666
$ws.Range("L28").Value = -5

$excel.CutCopyMode = 0
